# SignalR for full warehousedashboard
# Populate the receiving template header row with the real column names
# (Supplier, ReceiptDate, Note, ProductName, Uom, Quantity, UnitPrice),
# replacing the old single placeholder string, and size the columns to
# comfortably fit that content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Supplier"
$ws.Range("B1").Value = "ReceiptDate"
$ws.Range("C1").Value = "Note"
$ws.Range("D1").Value = "ProductName"
$ws.Range("E1").Value = "Uom"
$ws.Range("F1").Value = "Quantity"
$ws.Range("G1").Value = "UnitPrice"

$ws.Columns.Item(1).ColumnWidth = 23.5
$ws.Columns.Item(2).ColumnWidth = 18.5
$ws.Columns.Item(3).ColumnWidth = 48.333333333333336
$ws.Columns.Item(4).ColumnWidth = 14.833333333333334
$ws.Columns.Item(5).ColumnWidth = 17.666666666666668
$ws.Columns.Item(6).ColumnWidth = 11.166666666666666
$ws.Columns.Item(7).ColumnWidth = 14.666666666666666

$ws.Columns.Item(7).OutlineLevel = 6

$ws.Range("H13").Select()
